$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "VERSION"
$newSheet.Range("A6").Value = "File version"
$newSheet.Range("B6").Value = "vx.xx"
$newSheet.Activate()
$newSheet.Range("B6").Select()
Write-Host $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
